$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 80. This pushes the existing rows
# 80-85 down to 81-86, preserving their values and formatting.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with its data.
$row = 80
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 45106
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112031
$ws.Cells.Item($row, 7).Value = "Poroto verde"
$ws.Cells.Item($row, 8).Value = "Magnum"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 14000
$ws.Cells.Item($row, 12).Value = 15000
$ws.Cells.Item($row, 13).Value = 14500
$ws.Cells.Item($row, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 580
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
